$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 202.9375
$ws.Range("J12").Value = 148.5
$ws.Range("L12").Value = 148.5
$ws.Range("N12").Value = -488.5

$ws.Range("H62").Value = 10470
$ws.Range("I62").Value = 8696.25
$ws.Range("J62").Value = 17565
$ws.Range("K62").Value = 8696.25
$ws.Range("L62").Value = 17565
$ws.Range("M62").Value = -8072.25
$ws.Range("N62").Value = -18813

$ws.Range("H64").Value = 16045.15
$ws.Range("I64").Value = 20792.615
$ws.Range("J64").Value = 7228.4287
$ws.Range("K64").Value = 20792.615
$ws.Range("L64").Value = 7228.4287
$ws.Range("M64").Value = -20544.615
$ws.Range("N64").Value = -7724.4287

$ws.Range("H65").Value = 10470
$ws.Range("I65").Value = 8696.25
$ws.Range("J65").Value = 17565
$ws.Range("K65").Value = 43481.25
$ws.Range("L65").Value = 87825
$ws.Range("M65").Value = -40361.25
$ws.Range("N65").Value = -94065

$ws.Range("H67").Value = 16045.15
$ws.Range("I67").Value = 20792.615
$ws.Range("J67").Value = 7228.4287
$ws.Range("K67").Value = 20792.615
$ws.Range("L67").Value = 7228.4287
$ws.Range("M67").Value = -19934.615
$ws.Range("N67").Value = -8944.4287

$ws.Range("H74").Value = 9606.799999999999
$ws.Range("I74").Value = 9508.583000000001
$ws.Range("J74").Value = 9999.666999999999
$ws.Range("K74").Value = 9508.583000000001
$ws.Range("L74").Value = 9999.666999999999
$ws.Range("M74").Value = -8572.583000000001
$ws.Range("N74").Value = -11871.667

$ws.Range("H76").Value = 5260.8125
$ws.Range("I76").Value = 4247.6
$ws.Range("J76").Value = 5721.364
$ws.Range("K76").Value = 4247.6
$ws.Range("L76").Value = 5721.364
$ws.Range("M76").Value = -3932.6
$ws.Range("N76").Value = -6351.364

$ws.Range("H77").Value = 9606.799999999999
$ws.Range("I77").Value = 9508.583000000001
$ws.Range("J77").Value = 9999.666999999999
$ws.Range("K77").Value = 47542.915
$ws.Range("L77").Value = 49998.335
$ws.Range("M77").Value = -42862.915
$ws.Range("N77").Value = -59358.335

$ws.Range("H79").Value = 5260.8125
$ws.Range("I79").Value = 4247.6
$ws.Range("J79").Value = 5721.364
$ws.Range("K79").Value = 4247.6
$ws.Range("L79").Value = 5721.364
$ws.Range("M79").Value = -3155.6
$ws.Range("N79").Value = -7905.364

$ws.Range("H107").Value = 1197.9
$ws.Range("I107").Value = 1071.6666
$ws.Range("K107").Value = 1071.6666
$ws.Range("M107").Value = 848.3334

$ws.Range("H113").Value = 3850
$ws.Range("I113").Value = 3850
$ws.Range("K113").Value = 3850
$ws.Range("M113").Value = -596

$ws.Range("H132").Value = 15684.306
$ws.Range("I132").Value = 15684.306
$ws.Range("K132").Value = 47052.91800000001
$ws.Range("M132").Value = -44522.91800000001

$ws.Range("H137").Value = 17122.8
$ws.Range("I137").Value = 23040.072
$ws.Range("J137").Value = 3315.8333
$ws.Range("K137").Value = 69120.216
$ws.Range("L137").Value = 9947.499899999999
$ws.Range("M137").Value = -66570.216
$ws.Range("N137").Value = -15047.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26661.244
$ws.Range("I32").Value = 29429.19
$ws.Range("J32").Value = 1057.75
$ws.Range("K32").Value = 29429.19
$ws.Range("L32").Value = 1057.75
$ws.Range("M32").Value = -29142.19
$ws.Range("N32").Value = -1631.75

$ws.Range("H34").Value = 255000
$ws.Range("I34").Value = 500000
$ws.Range("K34").Value = 500000
$ws.Range("M34").Value = -499729

$ws.Range("H61").Value = 8537.267
$ws.Range("I61").Value = 1306.4445
$ws.Range("J61").Value = 19383.5
$ws.Range("K61").Value = 1306.4445
$ws.Range("L61").Value = 19383.5
$ws.Range("M61").Value = -1094.4445
$ws.Range("N61").Value = -19807.5

$ws.Range("H74").Value = 292632
$ws.Range("I74").Value = 463018.94
$ws.Range("J74").Value = 15753.25
$ws.Range("K74").Value = 463018.94
$ws.Range("L74").Value = 15753.25
$ws.Range("M74").Value = -462144.94
$ws.Range("N74").Value = -17501.25

$ws.Range("H77").Value = 292632
$ws.Range("I77").Value = 463018.94
$ws.Range("J77").Value = 15753.25
$ws.Range("K77").Value = 2315094.7
$ws.Range("L77").Value = 78766.25
$ws.Range("M77").Value = -2310726.7
$ws.Range("N77").Value = -87502.25

$ws.Range("H122").Value = 3348.926
$ws.Range("I122").Value = 3317.7778
$ws.Range("J122").Value = 3411.2222
$ws.Range("K122").Value = 9953.3334
$ws.Range("L122").Value = 10233.6666
$ws.Range("M122").Value = -7503.3334
$ws.Range("N122").Value = -15133.6666

$ws.Range("H132").Value = 1587.381
$ws.Range("I132").Value = 1263.0555
$ws.Range("J132").Value = 3533.3333
$ws.Range("K132").Value = 3789.1665
$ws.Range("L132").Value = 10599.9999
$ws.Range("M132").Value = -1259.1665
$ws.Range("N132").Value = -15659.9999

$ws.Range("H136").Value = 8537.267
$ws.Range("I136").Value = 1306.4445
$ws.Range("J136").Value = 19383.5
$ws.Range("K136").Value = 3919.3335
$ws.Range("L136").Value = 58150.5
$ws.Range("M136").Value = -1369.3335
$ws.Range("N136").Value = -63250.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1824.3334
$ws.Range("I86").Value = 1539.2
$ws.Range("K86").Value = 1539.2
$ws.Range("M86").Value = -416.2

$ws.Range("H89").Value = 1824.3334
$ws.Range("I89").Value = 1539.2
$ws.Range("K89").Value = 7696
$ws.Range("M89").Value = -2080

$ws.Range("H94").Value = 7144740
$ws.Range("J94").Value = 20002236
$ws.Range("L94").Value = 20002236
$ws.Range("N94").Value = -20003138

$ws.Range("H105").Value = 1743.4706
$ws.Range("I105").Value = 1227.4286
$ws.Range("J105").Value = 4151.6665
$ws.Range("K105").Value = 1227.4286
$ws.Range("L105").Value = 4151.6665
$ws.Range("M105").Value = 519.5714
$ws.Range("N105").Value = -7645.6665

$ws.Range("H107").Value = 4078.9565
$ws.Range("I107").Value = 4183.533
$ws.Range("J107").Value = 3882.875
$ws.Range("K107").Value = 4183.533
$ws.Range("L107").Value = 3882.875
$ws.Range("M107").Value = -2263.533
$ws.Range("N107").Value = -7722.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2873
$ws.Range("I105").Value = 1153.8334
$ws.Range("J105").Value = 4346.5713
$ws.Range("K105").Value = 1153.8334
$ws.Range("L105").Value = 4346.5713
$ws.Range("M105").Value = 593.1666
$ws.Range("N105").Value = -7840.5713

$ws.Range("H107").Value = 723.2727
$ws.Range("I107").Value = 657.2143
$ws.Range("J107").Value = 838.875
$ws.Range("K107").Value = 657.2143
$ws.Range("L107").Value = 838.875
$ws.Range("M107").Value = 1262.7857
$ws.Range("N107").Value = -4678.875

$ws.Range("H122").Value = 2389.182
$ws.Range("J122").Value = 2804.6667
$ws.Range("L122").Value = 8414.000100000001
$ws.Range("N122").Value = -13314.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 802.75
$ws.Range("J5").Value = 699
$ws.Range("L5").Value = 2097
$ws.Range("N5").Value = -2321

$ws.Range("H132").Value = 1593.6666
$ws.Range("I132").Value = 639.3333
$ws.Range("J132").Value = 2070.8333
$ws.Range("K132").Value = 5753.9997
$ws.Range("L132").Value = 18637.4997
$ws.Range("M132").Value = -3223.9997
$ws.Range("N132").Value = -23697.4997

$ws.Range("H135").Value = 802.75
$ws.Range("J135").Value = 699
$ws.Range("L135").Value = 6291
$ws.Range("N135").Value = -11361

$ws.Range("H136").Value = 2164
$ws.Range("I136").Value = 2164
$ws.Range("K136").Value = 6492
$ws.Range("M136").Value = -1392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5185.6787
$ws.Range("J80").Value = 8301.111000000001
$ws.Range("L80").Value = 8301.111000000001
$ws.Range("N80").Value = -10297.111

$ws.Range("H83").Value = 5185.6787
$ws.Range("J83").Value = 8301.111000000001
$ws.Range("L83").Value = 41505.55500000001
$ws.Range("N83").Value = -51489.55500000001

$ws.Range("H102").Value = 2580.8
$ws.Range("I102").Value = 2312
$ws.Range("K102").Value = 2312
$ws.Range("M102").Value = -690

$ws.Range("H113").Value = 2548.1538
$ws.Range("I113").Value = 1411.3
$ws.Range("K113").Value = 1411.3
$ws.Range("M113").Value = 758.7

$ws.Range("H122").Value = 41670320
$ws.Range("I122").Value = 3480.4443
$ws.Range("K122").Value = 10441.3329
$ws.Range("M122").Value = -7991.332900000001

$ws.Range("H126").Value = 1642.0667
$ws.Range("I126").Value = 1642.0667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4926.2001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2456.2001
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2779.6924
$ws.Range("I132").Value = 2677.9167
$ws.Range("J132").Value = 4001
$ws.Range("K132").Value = 8033.750100000001
$ws.Range("L132").Value = 12003
$ws.Range("M132").Value = -5503.750100000001
$ws.Range("N132").Value = -17063

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1889.1538
$ws.Range("I40").Value = 1833.6666
$ws.Range("K40").Value = 1833.6666
$ws.Range("M40").Value = -1697.6666

$ws.Range("H132").Value = 5263.6
$ws.Range("I132").Value = 5087.1816
$ws.Range("J132").Value = 5748.75
$ws.Range("K132").Value = 15261.5448
$ws.Range("L132").Value = 17246.25
$ws.Range("M132").Value = -12731.5448
$ws.Range("N132").Value = -22306.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 53994.32
$ws.Range("I122").Value = 62318.625
$ws.Range("K122").Value = 186955.875
$ws.Range("M122").Value = -184505.875

$ws.Range("H136").Value = 19058.484
$ws.Range("I136").Value = 24876.209
$ws.Range("K136").Value = 74628.62699999999
$ws.Range("M136").Value = -72078.62699999999
